$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prefix")
$ws.Range("A1").Value = "hello"
